# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 45243 (2023-11-13) to 45244 (2023-11-14).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    if ($cell.Value2() -eq 45243) {
        $cell.Value2 = 45244
    }
}
